$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Logs": append a new row (row 8) with the 7th test-mail entry
# -----------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Ik ben niet tevreden met mijn bestelling."
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #7: Ik ben niet tevreden met mijn bestelling."
$logs.Range("D8").Value = "Retour / Terugbetaling"

$antwoord8 = "Beste klant,`n" +
  "Bedankt voor uw bericht. Wat vervelend om te horen dat u niet tevreden bent met uw bestelling. Om u beter van dienst te kunnen zijn, zou ik graag meer details willen weten over wat er precies niet naar wens is gegaan. Kunt u mogelijk informatie geven over het specifieke product of de reden waarom u niet tevreden bent? Op die manier kunnen we het probleem verder onderzoeken en een passende oplossing bieden.`n" +
  "Alvast bedankt voor uw medewerking.`n" +
  "Met vriendelijke groet,`n" +
  "[Naam]  `n" +
  "E-mailassistent  `n" +
  "[Bedrijfsnaam]"
$logs.Range("E8").Value = $antwoord8

$logs.Range("F8").Value = "2025-07-22 12:25:32"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Ja"

# Extend the conditional-formatting ranges so the new row is covered too
# (D2:D7 -> D2:D8, G2:G7 -> G2:G8, H2:H7 -> H2:H8, I2:I7 -> I2:I8, J2:J7 -> J2:J8)
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
  $oldRange = $logs.Range($col + "2:" + $col + "7")
  $newRange = $logs.Range($col + "2:" + $col + "8")
  $fcs = $oldRange.FormatConditions
  for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($newRange)
  }
}

# -----------------------------------------------------------------
# Sheet "Dashboard": the summary counts change because the new Logs
# row bumps "Retour / Terugbetaling" to 3, which now outranks
# "Productinformatie" (2), so the two rows swap places.
# -----------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Retour / Terugbetaling"
$dash.Range("B2").Value = 3
$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 2
